$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled header string (MODEL_CONDITION -> MODELCONDITION)
$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION")

# The original column A (bold-styled row index values 1/8) was a leftover
# artifact; delete it so every other column shifts left by one (B->A,
# C->B, D->C, E->D, F->E) and the used range becomes A1:E3.
$ws.Columns.Item(1).Delete()
